$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AT4 buff: row 15 - Damage and Damage Upgraded values get a 5th tier added
$ws.Range("D15").Value = 325290255220195
$ws.Range("E15").Value = 375340305265230

# AT4 Rate of Fire / Explosion Radius text updated with the new 5th shot tier
$ws.Range("G15").Value = "1 Frame / 35,70,105,140,175 Up 50,85,120,155,190"

# Column G needs to widen (bestFit) to accommodate the longer text
$ws.Columns.Item(7).ColumnWidth = 43.85

# Selection state left on column H by the editor during cleanup
$ws.Columns.Item(8).Select()
